$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# Footer note text changed from "name, email address" to "Bla Bla Bla"
$ws.Range("A31").Value = "Bla Bla Bla"

# Labor hours amount (E18) changed from 1 to 125
$ws.Range("E18").Value = 125
